# Rename the existing "Games" sheet to "Games PS4" and add a new
# "Games SWITCH" sheet with the Nintendo Switch catalogue.

$wb = $excel.ActiveWorkbook

$gamesPs4 = $wb.Worksheets.Item("Games")
$gamesPs4.Name = "Games PS4"

$switch = $wb.Worksheets.Add()
$switch.Name = "Games SWITCH"
# Move the newly created sheet to the end of the tab strip, after "Games PS4".
$switch.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# Re-resolve the worksheet by name after the move so subsequent writes
# land on the correct physical sheet.
$switch = $wb.Worksheets.Item("Games SWITCH")

$data = @(
    @("Game", "Price"),
    @("JOGO POKÉMON SWORD NINTENDO SWITCH", "R$ 369,49"),
    @("JOGO MARIO KART 8 DELUXE NINTENDO SWITCH", "R$ 325,51"),
    @("JOGO POKÉMON LEGENDS: ARCEUS NINTENDO SWITCH", "R$ 443,43"),
    @("JOGO SUPER SMASH BROS ULTIMATE NINTENDO SWITCH", "R$ 369,49"),
    @("JOGO POKÉMON LETS GO EEVEE NINTENDO SWITCH", "R$ 369,49"),
    @("JOGO NEW SUPER MARIO BROS. U DELUXE NINTENDO SWITCH", "R$ 369,49"),
    @("JOGO POKÉMON SHINING PEARL NINTENDO SWITCH", "R$ 369,49"),
    @("JOGO THE LEGEND OF ZELDA: TEARS OF THE KINGDOM NINTENDO SWITCH", "R$ 325,51"),
    @("JOGO MINECRAFT NINTENDO SWITCH", "R$ 312,19"),
    @("JOGO ANIMAL CROSSING: NEW HORIZONS", "R$ 307,91"),
    @("JOGO SUPER MARIO ODYSSEY NINTENDO SWITCH", "R$ 307,91"),
    @("JOGO MARIO STRIKERS BATTLE LEAGUE NINTENDO SWITCH", "R$ 263,91")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 1
    $switch.Cells.Item($row, 1).Value = $data[$i][0]
    $switch.Cells.Item($row, 2).Value = $data[$i][1]
}
